$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.943.03'
$ws.Range("E2").Value = '  -0.04%  '
$ws.Range("D3").Value = '1.852.35'
$ws.Range("E3").Value = '  -0.94%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  +0.45%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.49'
$ws.Range("E5").Value = '  -0.20%  '
$ws.Range("E6").Value = '  +0.33%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5057'
$ws.Range("E7").Value = '  +1.30%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3801'
$ws.Range("E8").Value = '  -0.35%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08218'
$ws.Range("E9").Value = '  -8.05%  '
$ws.Range("B10").Value = 'OKB'
$ws.Range("C10").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.53'
$ws.Range("E10").Value = '  +0.24%  '
$ws.Range("B11").Value = 'Polygon'
$ws.Range("C11").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.104'
$ws.Range("E11").Value = '  -1.39%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.175'
$ws.Range("E12").Value = '  -2.25%  '
$ws.Range("D13").Value = '1.865.03'
$ws.Range("E13").Value = '  -0.21%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.39'
$ws.Range("E14").Value = '  -1.41%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.173'
$ws.Range("E15").Value = '  -0.93%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.005'
$ws.Range("E16").Value = '  +0.54%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001092'
$ws.Range("E17").Value = '  -0.85%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '90.23'
$ws.Range("E18").Value = '  -0.70%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06609'
$ws.Range("E19").Value = '  -0.45%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.64'
$ws.Range("E20").Value = '  -1.58%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.002'
$ws.Range("E21").Value = '  +0.24%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.989'
$ws.Range("E22").Value = '  -1.98%  '
$ws.Range("D23").Value = '27.964.06'
$ws.Range("E23").Value = '  -0.04%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.00'
$ws.Range("E24").Value = '  -4.39%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.243'
$ws.Range("E25").Value = '  -1.64%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.532'
$ws.Range("E26").Value = '  +0.66%  '
$ws.Range("D27").Value = '2.074.94'
$ws.Range("E27").Value = '  -0.09%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '157.60'
$ws.Range("E28").Value = '  -0.30%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '20.31'
$ws.Range("E29").Value = '  -1.92%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '123.90'
$ws.Range("E30").Value = '  -1.92%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.1054'
$ws.Range("E31").Value = '  -0.47%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.031'
$ws.Range("E32").Value = '  -2.37%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.577'
$ws.Range("E33").Value = '  +0.01%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.595'
$ws.Range("E34").Value = '  +0.15%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.535'
$ws.Range("E35").Value = '  +1.78%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06494'
$ws.Range("E36").Value = '  -0.61%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02398'
$ws.Range("E37").Value = '  +0.12%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2158'
$ws.Range("E38").Value = '  -1.62%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.195'
$ws.Range("E39").Value = '  -0.53%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.234'
$ws.Range("E40").Value = '  -5.85%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6354'
$ws.Range("E41").Value = '  -0.31%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.12'
$ws.Range("E42").Value = '  -4.66%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.836'
$ws.Range("E43").Value = '  -1.36%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6024'
$ws.Range("E44").Value = '  +0.29%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.10'
$ws.Range("E45").Value = '  -1.21%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.280'
$ws.Range("E46").Value = '  +0.09%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.650'
$ws.Range("E47").Value = '  -0.60%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.969'
$ws.Range("E48").Value = '  -0.54%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.201'
$ws.Range("E49").Value = '  -2.03%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '120.42'
$ws.Range("E50").Value = '  -0.95%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '78.62'
$ws.Range("E51").Value = '  +0.62%  '
